# trafo_id -> gridnode_id refactor
# The "electrolysers" sheet has a header row (row 1) listing field names.
# Column J currently holds "trafo_id"; rename it to "gridnode_id".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("electrolysers")

$ws.Range("J1").Value = "gridnode_id"

# Update the active cell selection recorded in the sheet view.
$ws.Range("G7").Select()
